$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth
$ws.Activate()
$ws.Range("K14").Select()
